$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 2 and Row 3 had their Id/Antal/Ost/Nord/Publik-kommentar values swapped.
# Column A (Id)
$ws.Range("A2").Value = 111863045
$ws.Range("A3").Value = 111863073

# Column I (Antal) - stored as text in the source data, so force text type
$ws.Range("I2").NumberFormat = "General"
$ws.Range("I2").Value = "'11"
$ws.Range("I3").NumberFormat = "General"
$ws.Range("I3").Value = "'2"

# Column Q (Ost)
$ws.Range("Q2").Value = 655234
$ws.Range("Q3").Value = 655228

# Column R (Nord)
$ws.Range("R2").Value = 6634889
$ws.Range("R3").Value = 6634879

# Column AC (Publik kommentar)
$ws.Range("AC2").Value = "Under gran i svacka"
$ws.Range("AC3").Value = ""
